$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35
$ws.Range("N35").Value = 13

# Row 38
$ws.Range("H38").Value = 3
$ws.Range("I38").Value = 3.1
$ws.Range("L38").Value = 4
$ws.Range("O38").Value = 1.53
$ws.Range("P38").Value = 2.5
$ws.Range("U38").Value = 2.1
$ws.Range("V38").Value = 1.67
$ws.Range("Y38").Value = 11
$ws.Range("Z38").Value = 26
$ws.Range("AV38").Value = 81
$ws.Range("BB38").Value = 351

# Row 39
$ws.Range("Q39").Value = 2.1
$ws.Range("R39").Value = 1.7

# Row 41
$ws.Range("G41").Value = 3
$ws.Range("I41").Value = 2.63
$ws.Range("J41").Value = 3.75
$ws.Range("L41").Value = 3.5
$ws.Range("W41").Value = 7
$ws.Range("X41").Value = 13
$ws.Range("Y41").Value = 12
$ws.Range("Z41").Value = 29
$ws.Range("AA41").Value = 29
$ws.Range("AG41").Value = 1250
$ws.Range("AH41").Value = 6.5
$ws.Range("AI41").Value = 11
$ws.Range("AK41").Value = 26
$ws.Range("AL41").Value = 26
$ws.Range("AN41").Value = 4.75
$ws.Range("AW41").Value = 4.5
$ws.Range("AX41").Value = 17
$ws.Range("AY41").Value = 34
$ws.Range("BA41").Value = 101

# Row 87
$ws.Range("G87").Value = 2.32
$ws.Range("H87").Value = 3.1
$ws.Range("J87").Value = 2.77
$ws.Range("L87").Value = 3.3
$ws.Range("M87").Value = 1.02
$ws.Range("N87").Value = 13
$ws.Range("P87").Value = 4.5
$ws.Range("Q87").Value = 1.53
$ws.Range("R87").Value = 2.18
$ws.Range("S87").Value = 1.3
$ws.Range("T87").Value = 3.4
$ws.Range("U87").Value = 1.54
$ws.Range("V87").Value = 2.38
$ws.Range("W87").Value = 10
$ws.Range("X87").Value = 13
$ws.Range("Y87").Value = 7.7
$ws.Range("Z87").Value = 23
$ws.Range("AA87").Value = 13.5
$ws.Range("AB87").Value = 15
$ws.Range("AG87").Value = 110
$ws.Range("AH87").Value = 10.25
$ws.Range("AI87").Value = 14.5
$ws.Range("AK87").Value = 29
$ws.Range("AL87").Value = 17.5
$ws.Range("AM87").Value = 18
$ws.Range("AP87").Value = 15
$ws.Range("AR87").Value = 55
$ws.Range("AS87").Value = 120
$ws.Range("AX87").Value = 15
$ws.Range("AY87").Value = 18
$ws.Range("AZ87").Value = 65
$ws.Range("BA87").Value = 75
$ws.Range("BB87").Value = 175
